# Fix the StylishResume.docx template bug: the lone "{Projects}" run needs
# to be split into four separate runs -- "{", "{", "Projects}", "}" -- so the
# downstream doc_from_template.py parser (which walks runs looking for the
# "{{ ... }}" placeholder pattern split across run boundaries) sees the
# expected token shape instead of a single merged run.

$d = $word.ActiveDocument

# Locate the exact "{Projects}" run (the lone-brace placeholder), not the
# "{% if Projects %}" / "{{Education}}" style text elsewhere in the doc.
$findRange = $d.Content
$found = $findRange.Find.Execute("{Projects}", $true, $true, $false, $false, `
                                  $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the '{Projects}' placeholder run to fix."
}
if ($findRange.Text -ne "{Projects}") {
    throw "Unexpected match text: [$($findRange.Text)]"
}

# This targets the known "{Projects}" placeholder paragraph in
# StylishResume.docx; keep its original w14:paraId/w14:textId/rsid*
# identity attributes on the rewritten paragraph instead of letting
# InsertXML mint fresh ones.
$pAttrs = ' w14:paraId="47B23FDB" w14:textId="1B63CA9E" w:rsidR="00FB0C4B"' + `
          ' w:rsidRDefault="0005561A" w:rsidP="0005561A"'

# Replace the single "{Projects}" run with four runs: "{", "{", "Projects}",
# "}" -- same visible text, same (default) formatting, just split apart.
$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" ' + `
       'xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"' + $pAttrs + '>' + `
         '<w:r><w:t>{</w:t></w:r>' + `
         '<w:r><w:t>{</w:t></w:r>' + `
         '<w:r><w:t>Projects}</w:t></w:r>' + `
         '<w:r><w:t>}</w:t></w:r>' + `
       '</w:p>'

$null = $findRange.InsertXML($xml)

Write-Output "Replaced '{Projects}' run with 4-run split."
